$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-66 down to 19-67.
$ws.Rows('18:18').Insert()

# Populate the newly inserted row 18 with the new weekly price entry.
$ws.Cells.Item(18, 1).Value = 2
$ws.Cells.Item(18, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(18, 3).Value = 'Coquimbo'
$ws.Cells.Item(18, 4).Value = 44581
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = 100112030
$ws.Cells.Item(18, 7).Value = 'Poroto granado'
$ws.Cells.Item(18, 8).Value = 'Sin especificar'
$ws.Cells.Item(18, 9).Value = 'Primera'
$ws.Cells.Item(18, 10).Value = 600
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 15000
$ws.Cells.Item(18, 13).Value = 14000
$ws.Cells.Item(18, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(18, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(18, 16).Value = 933
$ws.Cells.Item(18, 17).Value = 15
$ws.Cells.Item(18, 18).Value = 'Hortaliza'
